$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.64
$ws.Range("G2").Value = 2.08
$ws.Range("I2").Value = 7.4
$ws.Range("K2").Value = 6.4
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 1.92
$ws.Range("N3").Value = 2.06
$ws.Range("P3").Value = 2.06
$ws.Range("O6").Value = 1.43
$ws.Range("Y6").Value = 10.5
$ws.Range("AA8").Value = 980
$ws.Range("AB8").Value = 13.5
$ws.Range("AC8").Value = 8.800000000000001
$ws.Range("AD8").Value = 13
$ws.Range("AE8").Value = 980
$ws.Range("AF8").Value = 21
$ws.Range("AG8").Value = 13.5
$ws.Range("AH8").Value = 17
$ws.Range("AI8").Value = 980
$ws.Range("AJ8").Value = 980
$ws.Range("AK8").Value = 980
$ws.Range("AL8").Value = 980
$ws.Range("AM8").Value = 85
$ws.Range("AN8").Value = 25
$ws.Range("AO8").Value = 980
$ws.Range("F8").Value = 2.62
$ws.Range("G8").Value = 3.05
$ws.Range("H8").Value = 2.54
$ws.Range("I8").Value = 2.9
$ws.Range("J8").Value = 3.4
$ws.Range("K8").Value = 3.95
$ws.Range("N8").Value = 3.9
$ws.Range("O8").Value = 1.27
$ws.Range("P8").Value = 2.02
$ws.Range("Q8").Value = 1.81
$ws.Range("R8").Value = 1.39
$ws.Range("S8").Value = 3
$ws.Range("T8").Value = 1.67
$ws.Range("U8").Value = 2.22
$ws.Range("V8").Value = 1.52
$ws.Range("W8").Value = 1.49
$ws.Range("X8").Value = 19
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 20
$ws.Range("L9").Value = 1.18
$ws.Range("R9").Value = 1.81
$ws.Range("U9").Value = 1.81
$ws.Range("AA10").Value = 360
$ws.Range("AB10").Value = 14
$ws.Range("AC10").Value = 16
$ws.Range("AD10").Value = 970
$ws.Range("AE10").Value = 150
$ws.Range("AF10").Value = 11.5
$ws.Range("AG10").Value = 12
$ws.Range("AH10").Value = 970
$ws.Range("AI10").Value = 120
$ws.Range("AJ10").Value = 13
$ws.Range("AK10").Value = 14.5
$ws.Range("AL10").Value = 970
$ws.Range("AM10").Value = 120
$ws.Range("AN10").Value = 4.1
$ws.Range("AO10").Value = 140
$ws.Range("F10").Value = 1.33
$ws.Range("G10").Value = 1.34
$ws.Range("L10").Value = 1.22
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 1.13
$ws.Range("P10").Value = 2.94
$ws.Range("Q10").Value = 1.35
$ws.Range("R10").Value = 1.8
$ws.Range("S10").Value = 1.86
$ws.Range("T10").Value = 1.74
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.1
$ws.Range("X10").Value = 970
$ws.Range("Y10").Value = 970
$ws.Range("Z10").Value = 120
$ws.Range("AH11").Value = 15
$ws.Range("AN11").Value = 23
$ws.Range("H11").Value = 2.08
$ws.Range("AB12").Value = 19.5
$ws.Range("AJ12").Value = 60
$ws.Range("AN12").Value = 19.5
$ws.Range("AO12").Value = 11
$ws.Range("F12").Value = 3.2
$ws.Range("G12").Value = 3.25
$ws.Range("H12").Value = 2.2
$ws.Range("N12").Value = 6
$ws.Range("S12").Value = 2.46
$ws.Range("U12").Value = 2.72
$ws.Range("V12").Value = 1.81
$ws.Range("W12").Value = 1.44
$ws.Range("AC13").Value = 7.8
$ws.Range("AK13").Value = 34
$ws.Range("G13").Value = 2.9
$ws.Range("O13").Value = 1.38
$ws.Range("Q13").Value = 2.1
$ws.Range("G14").Value = 2.14
$ws.Range("J14").Value = 2.84
$ws.Range("L14").Value = 1.6
$ws.Range("N14").Value = 2.12
$ws.Range("S14").Value = 8
$ws.Range("W14").Value = 1.87
$ws.Range("AB15").Value = 12.5
$ws.Range("AC15").Value = 8.800000000000001
$ws.Range("AE15").Value = 38
$ws.Range("AF15").Value = 19
$ws.Range("F15").Value = 2.34
$ws.Range("G15").Value = 2.54
$ws.Range("H15").Value = 3
$ws.Range("J15").Value = 3.5
$ws.Range("K15").Value = 3.85
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 4
$ws.Range("O15").Value = 1.26
$ws.Range("P15").Value = 2.04
$ws.Range("Q15").Value = 1.78
$ws.Range("R15").Value = 1.41
$ws.Range("S15").Value = 2.94
$ws.Range("T15").Value = 1.64
$ws.Range("U15").Value = 2.24
$ws.Range("V15").Value = 1.45
$ws.Range("W15").Value = 1.65
$ws.Range("X15").Value = 17.5
$ws.Range("AA16").Value = 240
$ws.Range("AB16").Value = 12.5
$ws.Range("AC16").Value = 13.5
$ws.Range("AD16").Value = 34
$ws.Range("AE16").Value = 120
$ws.Range("AF16").Value = 12
$ws.Range("AG16").Value = 12
$ws.Range("AH16").Value = 26
$ws.Range("AI16").Value = 95
$ws.Range("AJ16").Value = 15.5
$ws.Range("AK16").Value = 17
$ws.Range("AL16").Value = 36
$ws.Range("AM16").Value = 120
$ws.Range("AN16").Value = 6.8
$ws.Range("AO16").Value = 120
$ws.Range("F16").Value = 1.43
$ws.Range("G16").Value = 1.54
$ws.Range("H16").Value = 6.8
$ws.Range("I16").Value = 9.199999999999999
$ws.Range("J16").Value = 4.3
$ws.Range("K16").Value = 5.6
$ws.Range("L16").Value = 1.26
$ws.Range("N16").Value = 5
$ws.Range("O16").Value = 1.19
$ws.Range("P16").Value = 2.38
$ws.Range("Q16").Value = 1.57
$ws.Range("R16").Value = 1.55
$ws.Range("S16").Value = 2.44
$ws.Range("T16").Value = 1.77
$ws.Range("U16").Value = 2.04
$ws.Range("W16").Value = 2.84
$ws.Range("X16").Value = 28
$ws.Range("Y16").Value = 34
$ws.Range("Z16").Value = 80
$ws.Range("AA17").Value = 95
$ws.Range("AE17").Value = 55
$ws.Range("AM17").Value = 95
$ws.Range("AO17").Value = 55
$ws.Range("F17").Value = 1.88
$ws.Range("G17").Value = 2.1
$ws.Range("H17").Value = 3.7
$ws.Range("I17").Value = 4.6
$ws.Range("T17").Value = 1.68
$ws.Range("U17").Value = 2.18
$ws.Range("V17").Value = 1.28
$ws.Range("W17").Value = 1.91
$ws.Range("L18").Value = 1.24
$ws.Range("AA19").Value = 13.5
$ws.Range("AB19").Value = 46
$ws.Range("AC19").Value = 14.5
$ws.Range("AF19").Value = 90
$ws.Range("AG19").Value = 32
$ws.Range("AH19").Value = 21
$ws.Range("AJ19").Value = 260
$ws.Range("AK19").Value = 95
$ws.Range("AL19").Value = 75
$ws.Range("AO19").Value = 3.85
$ws.Range("F19").Value = 8.4
$ws.Range("G19").Value = 8.800000000000001
$ws.Range("H19").Value = 1.38
$ws.Range("I19").Value = 1.39
$ws.Range("J19").Value = 6.2
$ws.Range("K19").Value = 6.4
$ws.Range("L19").Value = 1.23
$ws.Range("N19").Value = 8.199999999999999
$ws.Range("O19").Value = 1.13
$ws.Range("P19").Value = 3.4
$ws.Range("Q19").Value = 1.4
$ws.Range("R19").Value = 1.96
$ws.Range("S19").Value = 1.99
$ws.Range("T19").Value = 1.66
$ws.Range("V19").Value = 3.55
$ws.Range("W19").Value = 1.12
$ws.Range("Y19").Value = 15
$ws.Range("N20").Value = 5.1
$ws.Range("Q20").Value = 1.69
$ws.Range("S20").Value = 2.72
$ws.Range("V20").Value = 1.92
$ws.Range("AA21").Value = 100
$ws.Range("AB21").Value = 8.6
$ws.Range("AD21").Value = 17.5
$ws.Range("AE21").Value = 60
$ws.Range("AG21").Value = 10.5
$ws.Range("AL21").Value = 40
$ws.Range("AN21").Value = 15
$ws.Range("F21").Value = 1.96
$ws.Range("G21").Value = 1.99
$ws.Range("J21").Value = 3.65
$ws.Range("K21").Value = 3.7
$ws.Range("L21").Value = 1.45
$ws.Range("N21").Value = 3.6
$ws.Range("P21").Value = 1.88
$ws.Range("T21").Value = 1.91
$ws.Range("W21").Value = 2
$ws.Range("Y21").Value = 15
$ws.Range("AI22").Value = 60
$ws.Range("G22").Value = 1.5
$ws.Range("I22").Value = 6.8
$ws.Range("J22").Value = 5.5
$ws.Range("K22").Value = 5.6
$ws.Range("W22").Value = 3
$ws.Range("AF23").Value = 9.199999999999999
$ws.Range("R23").Value = 1.98
$ws.Range("T23").Value = 2.18
$ws.Range("R24").Value = 1.35
$ws.Range("AB25").Value = 18
$ws.Range("AC25").Value = 29
$ws.Range("AD25").Value = 1000
$ws.Range("AF25").Value = 11
$ws.Range("AG25").Value = 14.5
$ws.Range("AH25").Value = 140
$ws.Range("AK25").Value = 13
$ws.Range("AL25").Value = 75
$ws.Range("AN25").Value = 2.36
$ws.Range("F25").Value = 1.13
$ws.Range("G25").Value = 1.15
$ws.Range("H25").Value = 24
$ws.Range("I25").Value = 26
$ws.Range("J25").Value = 12
$ws.Range("K25").Value = 13
$ws.Range("N25").Value = 11.5
$ws.Range("P25").Value = 4.4
$ws.Range("R25").Value = 2.38
$ws.Range("S25").Value = 1.69
$ws.Range("T25").Value = 2
$ws.Range("U25").Value = 1.9
$ws.Range("V25").Value = 1.04
$ws.Range("W25").Value = 8
